# locker and seatNumber bug fixed e.g in one locker only one student will be...
# Replace the sample data rows (2 and 3) on the "Students" sheet with the
# corrected values. All source columns are plain text (shared strings) in
# the workbook, including columns that look numeric/date-like (registration
# numbers, dates, amounts, seat/locker numbers), so we force a Text number
# format before writing each value and then clear the formatting again so
# the cells fall back to the workbook's default style - this keeps Excel
# from silently re-typing "1" as a number or "2025-03-07" as a date serial.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddress, $text) {
    $rng = $ws.Range($rangeAddress)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

# Row 2
Set-TextValue "A2" "1"
Set-TextValue "B2" "2025-03-07"
Set-TextValue "C2" "nikhil"
Set-TextValue "D2" "kjkhuhuj"
Set-TextValue "E2" "nghkk"
Set-TextValue "F2" "8651993323"
Set-TextValue "G2" "06:00-10:00, 22:00-06:00"
Set-TextValue "H2" "2"
Set-TextValue "I2" "1"
Set-TextValue "J2" "350.00"
Set-TextValue "K2" "50.00"
Set-TextValue "L2" "1"
Set-TextValue "M2" "2025-03-08"

# Row 3
Set-TextValue "A3" "2"
Set-TextValue "B3" "2025-02-05"
Set-TextValue "C3" "satyam"
Set-TextValue "D3" "satyam"
Set-TextValue "E3" "ramkrishnanagar"
Set-TextValue "F3" "7250585057"
Set-TextValue "G3" "10:00-14:00, 14:00-18:00"
Set-TextValue "H3" "2"
Set-TextValue "I3" "1"
Set-TextValue "J3" "400.00"
Set-TextValue "K3" "10.00"
Set-TextValue "L3" "1"
Set-TextValue "M3" "2025-03-12"

Write-Host "Applied locker/seatNumber fix sample data to rows 2-3"
